# Insert a new data row into the weekly "Feria Lagunitas de Puerto Montt -
# Zapallo" price table. The new record is inserted above the current row
# 263, shifting the existing rows 263-329 down to 264-330.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 263, pushing everything below it (263:329) down by one.
$ws.Rows("263:263").Insert()

# Populate the newly inserted row 263 with the new weekly record.
$ws.Cells.Item(263, 1).Value = 4
$ws.Cells.Item(263, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(263, 3).Value = "Los Lagos"
$ws.Cells.Item(263, 4).Value = 44754
$ws.Cells.Item(263, 5).Value = 10
$ws.Cells.Item(263, 6).Value = 100112045
$ws.Cells.Item(263, 7).Value = "Zapallo"
$ws.Cells.Item(263, 8).Value = "Paine"
$ws.Cells.Item(263, 9).Value = "1a (guarda)"
$ws.Cells.Item(263, 10).Value = 800
$ws.Cells.Item(263, 11).Value = 500
$ws.Cells.Item(263, 12).Value = 500
$ws.Cells.Item(263, 13).Value = 500
$ws.Cells.Item(263, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(263, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(263, 16).Value = 500
$ws.Cells.Item(263, 17).Value = 1
$ws.Cells.Item(263, 18).Value = "Hortaliza"
